# CANGUCU.xlsx update:
#  - Drop the "Desarquivamentos Pendentes" tab entirely.
#  - Rename "Paineis DARQ" -> "PAINEIS DARQ".
#  - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO".

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

$wsDesarquivamentos = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$wsDesarquivamentos.Delete()

$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"
